$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "safari degree added with so delay"
# Add a new student row (Mojtaba Safari Mohammad Abadi, row 15 already has the
# name in C15) with ID, Q1-Q4 scores and HW3/HW4 grades (HW3/HW4 delayed, so
# they get a "...per" (percentage-penalised) grade instead of a plain number).
# ---------------------------------------------------------------------------

# xlPasteFormats
$xlPasteFormats = -4122
# xlLineStyleNone / xlNone
$xlNone = -4142

# Template cell L3 already carries the "big Times-New-Roman, centered,
# shrink-to-fit" look (font id 2) with a real border that we strip off so the
# new cells end up borderless like the rest of the row.
$tpl = $ws.Range("L3")

# --- B15 / D15 / E15: plain (black) variant of the template style ---------
foreach ($addr in @("B15", "D15", "E15")) {
    $cell = $ws.Range($addr)
    $tpl.Copy()
    $cell.PasteSpecial($xlPasteFormats)
    $cell.Borders.LineStyle = $xlNone
}
$ws.Range("D15").Value = 97412301

# --- F15:K15 and F16:K16: red variant of the template style (delay marker) ---
$gradeRanges = @("F15:K15", "F16:K16")
foreach ($addr in $gradeRanges) {
    $cell = $ws.Range($addr)
    $tpl.Copy()
    $cell.PasteSpecial($xlPasteFormats)
    $cell.Borders.LineStyle = $xlNone
    $cell.Font.Color = 255
}

$ws.Range("F15").Value = 110
$ws.Range("G15").Value = 80
$ws.Range("H15").Value = 83
$ws.Range("I15").Value = 92
$ws.Range("K15").Value = "36.4 70per"
$ws.Range("J15").Value = "45 90per"

$excel.CutCopyMode = $false

# Move the active selection the way the saved workbook shows it.
$ws.Range("I16").Select()
